$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already contains a "tile" grid (the red/blue 16x18 tile) drawn
# with string borders at T20:AK37. Add a new ("white", shared string "111")
# tile right below it at T42:AK59, reusing the same border/interior layout.

# 1) Duplicate the existing tile block's formatting/structure into the new
#    location so the new tile matches the look (border of "000" cells
#    around an interior) of the other tiles on the sheet.
$ws.Range("T20:AK37").Copy()
$ws.Range("T42").PasteSpecial()

# 2) Re-label the interior of the new tile as "111" (the new shared string
#    for the new, third tile color/value) instead of the "100" used by the
#    tile that was copied from.
$ws.Range("U43:AJ58").Value = "111"

# 3) Update the selection to highlight the freshly added tile, matching
#    where the author's cursor ended up after adding it.
$ws.Range("T42:AK59").Select()
